$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain integers need an explicit
# text number format first, otherwise Excel COM auto-converts the
# assigned string into a numeric cell instead of keeping it textual.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "304"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "426"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "342"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "584"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "630"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "630"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "635"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "66"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "34"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "600"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "498"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "228"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "557"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "346"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "15"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "844"
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = "550"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "130"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "564"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "57"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "535"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "530"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "535"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "575"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "550"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "575"

# Remaining cells keep their default (unstyled) text representation;
# values containing spaces/commas are not auto-converted to numbers.
$ws.Range("C2").Value = "2 115"
$ws.Range("D2").Value = "1 790"
$ws.Range("E2").Value = "1 775"
$ws.Range("F2").Value = "1 760"
$ws.Range("G2").Value = "-1,68"
$ws.Range("C3").Value = "3 728"
$ws.Range("D3").Value = "5 245"
$ws.Range("E3").Value = "5 250"
$ws.Range("F3").Value = "5 245"
$ws.Range("G3").Value = "-0,10"
$ws.Range("D4").Value = "16 600"
$ws.Range("E4").Value = "16 600"
$ws.Range("F4").Value = "16 500"
$ws.Range("G4").Value = "-0,60"
$ws.Range("C5").Value = "2 242"
$ws.Range("D5").Value = "1 500"
$ws.Range("E5").Value = "1 445"
$ws.Range("F5").Value = "1 500"
$ws.Range("G5").Value = "7,14"
$ws.Range("D6").Value = "4 310"
$ws.Range("E6").Value = "4 350"
$ws.Range("F6").Value = "4 305"
$ws.Range("G6").Value = "0,12"
$ws.Range("C7").Value = "7 408"
$ws.Range("D7").Value = "3 520"
$ws.Range("E7").Value = "3 520"
$ws.Range("F7").Value = "3 510"
$ws.Range("G7").Value = "-0,28"
$ws.Range("C8").Value = "6 683"
$ws.Range("F8").Value = "6 050"
$ws.Range("G8").Value = "0,83"
$ws.Range("C9").Value = "15 104"
$ws.Range("D9").Value = "3 650"
$ws.Range("E9").Value = "3 695"
$ws.Range("F9").Value = "3 650"
$ws.Range("G9").Value = "-0,54"
$ws.Range("C10").Value = "5 002"
$ws.Range("D10").Value = "2 550"
$ws.Range("E10").Value = "2 510"
$ws.Range("F10").Value = "2 550"
$ws.Range("G10").Value = "1,59"
$ws.Range("C11").Value = "9 066"
$ws.Range("D11").Value = "4 650"
$ws.Range("E11").Value = "4 580"
$ws.Range("F11").Value = "4 645"
$ws.Range("G11").Value = "2,09"
$ws.Range("D12").Value = "1 290"
$ws.Range("F12").Value = "1 290"
$ws.Range("G12").Value = "-0,77"
$ws.Range("C13").Value = "3 530"
$ws.Range("D13").Value = "9 655"
$ws.Range("E13").Value = "9 900"
$ws.Range("F13").Value = "9 780"
$ws.Range("G13").Value = "-1,21"
$ws.Range("G14").Value = "-2,31"
$ws.Range("C15").Value = "3 983"
$ws.Range("D15").Value = "2 170"
$ws.Range("E15").Value = "2 170"
$ws.Range("F15").Value = "2 160"
$ws.Range("G15").Value = "0,47"
$ws.Range("C16").Value = "1 827"
$ws.Range("D16").Value = "11 200"
$ws.Range("F16").Value = "11 300"
$ws.Range("G16").Value = "0,89"
$ws.Range("C17").Value = "1 069 265"
$ws.Range("C18").Value = "4 189"
$ws.Range("D18").Value = "4 810"
$ws.Range("E18").Value = "4 850"
$ws.Range("F18").Value = "4 850"
$ws.Range("G18").Value = "0,83"
$ws.Range("E19").Value = "4 590"
$ws.Range("F19").Value = "4 550"
$ws.Range("G19").Value = "-0,98"
$ws.Range("G20").Value = "0,84"
$ws.Range("C21").Value = "1 101"
$ws.Range("D21").Value = "9 500"
$ws.Range("F21").Value = "9 595"
$ws.Range("G21").Value = "-0,05"
$ws.Range("C22").Value = "1 202"
$ws.Range("D22").Value = "12 300"
$ws.Range("E22").Value = "12 295"
$ws.Range("F22").Value = "12 475"
$ws.Range("G22").Value = "1,42"
$ws.Range("C23").Value = "4 673"
$ws.Range("D23").Value = "2 300"
$ws.Range("E23").Value = "2 260"
$ws.Range("F23").Value = "2 290"
$ws.Range("G23").Value = "-2,14"
$ws.Range("C24").Value = "8 942"
$ws.Range("D24").Value = "14 100"
$ws.Range("E24").Value = "14 100"
$ws.Range("F24").Value = "14 150"
$ws.Range("G24").Value = "0,35"
$ws.Range("D25").Value = "1 685"
$ws.Range("E25").Value = "1 625"
$ws.Range("F25").Value = "1 685"
$ws.Range("G25").Value = "4,01"
$ws.Range("C26").Value = "1 891"
$ws.Range("D26").Value = "8 795"
$ws.Range("E26").Value = "8 595"
$ws.Range("F26").Value = "8 800"
$ws.Range("G26").Value = "2,44"
$ws.Range("D27").Value = "2 680"
$ws.Range("E27").Value = "2 680"
$ws.Range("G27").Value = "0,00"
$ws.Range("C28").Value = "23 883"
$ws.Range("D28").Value = "1 210"
$ws.Range("E28").Value = "1 210"
$ws.Range("F28").Value = "1 210"
$ws.Range("G28").Value = "7,08"
$ws.Range("C29").Value = "14 053"
$ws.Range("D30").Value = "6 000"
$ws.Range("E30").Value = "6 000"
$ws.Range("F30").Value = "5 890"
$ws.Range("G30").Value = "-1,83"
$ws.Range("C31").Value = "1 147"
$ws.Range("F31").Value = "1 435"
$ws.Range("G31").Value = "-0,35"
$ws.Range("D33").Value = "25 400"
$ws.Range("E33").Value = "25 900"
$ws.Range("F33").Value = "25 400"
$ws.Range("G33").Value = "3,67"
$ws.Range("C34").Value = "10 457"
$ws.Range("D34").Value = "1 135"
$ws.Range("E34").Value = "1 140"
$ws.Range("F34").Value = "1 130"
$ws.Range("G34").Value = "-0,44"
$ws.Range("C35").Value = "8 728"
$ws.Range("E35").Value = "4 795"
$ws.Range("F35").Value = "4 790"
$ws.Range("G35").Value = "-0,10"
$ws.Range("E36").Value = "3 280"
$ws.Range("G37").Value = "1,85"
$ws.Range("C38").Value = "1 044"
$ws.Range("D38").Value = "18 000"
$ws.Range("E38").Value = "18 700"
$ws.Range("F38").Value = "18 000"
$ws.Range("G38").Value = "-3,23"
$ws.Range("D39").Value = "9 700"
$ws.Range("F39").Value = "9 700"
$ws.Range("G39").Value = "4,86"
$ws.Range("C40").Value = "23 574"
$ws.Range("D40").Value = "24 995"
$ws.Range("E40").Value = "24 750"
$ws.Range("F40").Value = "24 900"
$ws.Range("G40").Value = "-0,38"
$ws.Range("C41").Value = "14 826"
$ws.Range("D41").Value = "7 750"
$ws.Range("E41").Value = "7 750"
$ws.Range("G41").Value = "0,00"
$ws.Range("D42").Value = "6 510"
$ws.Range("E42").Value = "6 695"
$ws.Range("F42").Value = "6 600"
$ws.Range("G42").Value = "-1,42"
$ws.Range("G43").Value = "-2,73"
$ws.Range("C44").Value = "6 047"
$ws.Range("D44").Value = "20 795"
$ws.Range("E44").Value = "20 805"
$ws.Range("F44").Value = "20 795"
$ws.Range("G44").Value = "3,53"
$ws.Range("C45").Value = "49 595"
$ws.Range("D45").Value = "2 425"
$ws.Range("F45").Value = "2 425"
$ws.Range("G45").Value = "-6,55"
$ws.Range("C46").Value = "1 024"
$ws.Range("E46").Value = "2 400"
$ws.Range("C48").Value = "2 069"
$ws.Range("G48").Value = "4,55"
